$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 31251652
$ws.Range("I15").Value = 31251652
$ws.Range("K15").Value = 93754956
$ws.Range("M15").Value = -93754787

$ws.Range("H92").Value = 885.1667
$ws.Range("J92").Value = 5000
$ws.Range("L92").Value = 5000
$ws.Range("N92").Value = -7496

$ws.Range("H100").Value = 2204.5715
$ws.Range("J100").Value = 3593.4
$ws.Range("L100").Value = 3593.4
$ws.Range("N100").Value = -4675.4

$ws.Range("H116").Value = 14709950
$ws.Range("I116").Value = 35715640
$ws.Range("K116").Value = 35715640
$ws.Range("M116").Value = -35712198

$ws.Range("H123").Value = 79499.664
$ws.Range("J123").Value = 79499.664
$ws.Range("L123").Value = 79499.664
$ws.Range("N123").Value = -89299.664

$ws.Range("H125").Value = 38462252
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H132").Value = 1581.55
$ws.Range("I132").Value = 1646.7778
$ws.Range("K132").Value = 4940.3334
$ws.Range("M132").Value = -2410.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 11922696
$ws.Range("I97").Value = 1010
$ws.Range("J97").Value = 13909643
$ws.Range("K97").Value = 1010
$ws.Range("L97").Value = 13909643
$ws.Range("N97").Value = -13910635
$ws.Range("M97").Value = -514

$ws.Range("H122").Value = 9771.192999999999
$ws.Range("I122").Value = 10070.704
$ws.Range("K122").Value = 30212.112
$ws.Range("M122").Value = -27762.112

$ws.Range("H132").Value = 6754.3096
$ws.Range("I132").Value = 5067.7334
$ws.Range("J132").Value = 10970.75
$ws.Range("K132").Value = 15203.2002
$ws.Range("L132").Value = 32912.25
$ws.Range("M132").Value = -12673.2002
$ws.Range("N132").Value = -37972.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4506220.5
$ws.Range("I20").Value = 5377481.5
$ws.Range("K20").Value = 5377481.5
$ws.Range("M20").Value = -5377234.5

$ws.Range("H86").Value = 58827852
$ws.Range("I86").Value = 3125
$ws.Range("K86").Value = 3125
$ws.Range("M86").Value = -2002

$ws.Range("H89").Value = 58827852
$ws.Range("I89").Value = 3125
$ws.Range("K89").Value = 15625
$ws.Range("M89").Value = -10009

$ws.Range("H94").Value = 2112.7693
$ws.Range("I94").Value = 793.25
$ws.Range("K94").Value = 793.25
$ws.Range("M94").Value = -342.25

$ws.Range("H105").Value = 2773.7568
$ws.Range("I105").Value = 2383.5862
$ws.Range("K105").Value = 2383.5862
$ws.Range("M105").Value = -636.5862000000002

$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws.Range("H134").Value = 3673.3708
$ws.Range("I134").Value = 2203.1333
$ws.Range("K134").Value = 6609.3999
$ws.Range("M134").Value = -4074.3999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 127
$ws.Range("I19").Value = 127
$ws.Range("K19").Value = 127
$ws.Range("M19").Value = 43

$ws.Range("H24").Value = 127
$ws.Range("I24").Value = 127
$ws.Range("K24").Value = 127
$ws.Range("M24").Value = 43

$ws.Range("H31").Value = 6885.136
$ws.Range("I31").Value = 3002.389
$ws.Range("J31").Value = 11544.434
$ws.Range("K31").Value = 3002.389
$ws.Range("L31").Value = 11544.434
$ws.Range("M31").Value = -2707.389
$ws.Range("N31").Value = -12134.434

$ws.Range("H34").Value = 6885.136
$ws.Range("I34").Value = 3002.389
$ws.Range("J34").Value = 11544.434
$ws.Range("K34").Value = 3002.389
$ws.Range("L34").Value = 11544.434
$ws.Range("M34").Value = -2800.389
$ws.Range("N34").Value = -11948.434

$ws.Range("H86").Value = 3476346.2
$ws.Range("I86").Value = 4468160.5
$ws.Range("K86").Value = 4468160.5
$ws.Range("M86").Value = -4467037.5

$ws.Range("H89").Value = 3476346.2
$ws.Range("I89").Value = 4468160.5
$ws.Range("K89").Value = 22340802.5
$ws.Range("M89").Value = -22335186.5

$ws.Range("H115").Value = 59379
$ws.Range("J115").Value = 59379
$ws.Range("L115").Value = 59379
$ws.Range("N115").Value = -61729

$ws.Range("H119").Value = 95998
$ws.Range("J119").Value = 95998
$ws.Range("L119").Value = 95998
$ws.Range("N119").Value = -105674

$ws.Range("H122").Value = 1409.5
$ws.Range("I122").Value = 1181.7142
$ws.Range("K122").Value = 3545.1426
$ws.Range("M122").Value = -1095.1426

$ws.Range("H132").Value = 8824.333000000001
$ws.Range("I132").Value = 2378.4
$ws.Range("J132").Value = 13428.571
$ws.Range("K132").Value = 7135.200000000001
$ws.Range("L132").Value = 40285.713
$ws.Range("M132").Value = -4605.200000000001
$ws.Range("N132").Value = -45345.713

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 50005884
$ws.Range("I56").Value = 50005884
$ws.Range("K56").Value = 50005884
$ws.Range("M56").Value = -50005354

$ws.Range("H92").Value = 7693732.5
$ws.Range("J92").Value = 7693732.5
$ws.Range("L92").Value = 23081197.5
$ws.Range("N92").Value = -23083693.5

$ws.Range("H132").Value = 16624.889
$ws.Range("I132").Value = 11931.5
$ws.Range("J132").Value = 20379.6
$ws.Range("K132").Value = 107383.5
$ws.Range("L132").Value = 183416.4
$ws.Range("M132").Value = -104853.5
$ws.Range("N132").Value = -188476.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 75990
$ws.Range("I52").Value = 20000
$ws.Range("J52").Value = 89987.5
$ws.Range("K52").Value = 20000
$ws.Range("L52").Value = 89987.5
$ws.Range("M52").Value = -19741
$ws.Range("N52").Value = -90505.5

$ws.Range("H70").Value = 5804.5713
$ws.Range("I70").Value = 4629.853
$ws.Range("J70").Value = 8467.267
$ws.Range("K70").Value = 4629.853
$ws.Range("L70").Value = 8467.267
$ws.Range("M70").Value = -4359.853
$ws.Range("N70").Value = -9007.267

$ws.Range("H73").Value = 5804.5713
$ws.Range("I73").Value = 4629.853
$ws.Range("J73").Value = 8467.267
$ws.Range("K73").Value = 4629.853
$ws.Range("L73").Value = 8467.267
$ws.Range("M73").Value = -3693.853
$ws.Range("N73").Value = -10339.267

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H132").Value = 2401.077
$ws.Range("I132").Value = 2397.8823
$ws.Range("K132").Value = 7193.646900000001
$ws.Range("M132").Value = -4663.646900000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1397.4
$ws.Range("I16").Value = 1403.7391
$ws.Range("K16").Value = 1403.7391
$ws.Range("M16").Value = -1233.7391

$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws.Range("H46").Value = 13890438
$ws.Range("J46").Value = 18520284
$ws.Range("L46").Value = 18520284
$ws.Range("N46").Value = -18520660

$ws.Range("H80").Value = 49990
$ws.Range("J80").Value = 49990
$ws.Range("L80").Value = 49990
$ws.Range("N80").Value = -52236

$ws.Range("H83").Value = 49990
$ws.Range("J83").Value = 49990
$ws.Range("L83").Value = 149970
$ws.Range("N83").Value = -161202

$ws.Range("H93").Value = 6703.364
$ws.Range("I93").Value = 5500.143
$ws.Range("J93").Value = 8809
$ws.Range("K93").Value = 5500.143
$ws.Range("L93").Value = 8809
$ws.Range("M93").Value = -4252.143
$ws.Range("N93").Value = -11305

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 24998
$ws.Range("I15").Value = 24998
$ws.Range("K15").Value = 24998
$ws.Range("M15").Value = -24710

$ws.Range("H62").Value = 11980
$ws.Range("I62").Value = 11980
$ws.Range("K62").Value = 11980
$ws.Range("M62").Value = -11356

$ws.Range("H65").Value = 11980
$ws.Range("I65").Value = 11980
$ws.Range("K65").Value = 59900
$ws.Range("M65").Value = -56780

$ws.Range("H81").Value = 20013398
$ws.Range("I81").Value = 1662.8334
$ws.Range("J81").Value = 50031000
$ws.Range("K81").Value = 3325.6668
$ws.Range("L81").Value = 100062000
$ws.Range("M81").Value = -2264.6668
$ws.Range("N81").Value = -100064122

$ws.Range("H84").Value = 20013398
$ws.Range("I84").Value = 1662.8334
$ws.Range("J84").Value = 50031000
$ws.Range("K84").Value = 16628.334
$ws.Range("L84").Value = 500310000
$ws.Range("M84").Value = -11324.334
$ws.Range("N84").Value = -500320608

$ws.Range("H126").Value = 2979.4443
$ws.Range("I126").Value = 1799.8
$ws.Range("J126").Value = 4454
$ws.Range("K126").Value = 5399.4
$ws.Range("L126").Value = 13362
$ws.Range("M126").Value = -2929.4
$ws.Range("N126").Value = -18302

$ws.Range("H132").Value = 17880704
$ws.Range("I132").Value = 21745764
$ws.Range("K132").Value = 65237292
$ws.Range("M132").Value = -65234762
